# Tuntikirjanpito.xlsx - "small tweaks, mainly tuntikirjanpito"
# Continue the hour-bookkeeping log: move the "tunnit yht." (total hours)
# summary row down, and fill in a new batch of logged work entries for
# the backend/dev-environment work (rows 48-53), then re-total in row 60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear out the old summary row (currently row 50) ------------------
# It will be replaced by ordinary data rows; wipe its formatting/height too
# so it doesn't leave stale row metadata behind.
$ws.Rows("50:50").Clear()
$ws.Rows("50:50").AutoFit()

# --- 2. New log entries (rows 48-53) ---------------------------------------
# Row 48 starts a new day (7.12.2021 -> serial 44537), formatted like the
# other date cells (copy format only from A44, then set the value).
$ws.Range("A44").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A48").Value = 44537

$ws.Range("B48").Value = 1
$ws.Range("C48").Value = "backend perus error handling, logging tehty"
$ws.Range("D48").Value = "api"

$ws.Range("B49").Value = 1.5
$ws.Range("C49").Value = "backend pari scheemaa, perusasetuksia, user ja favoriteDate"
$ws.Range("D49").Value = "api"

$ws.Range("B50").Value = 1
$ws.Range("C50").Value = "usersRouter, error handling "
$ws.Range("D50").Value = "api"

$ws.Range("B51").Value = 1
$ws.Range("C51").Value = "asynchandler, usersRouter testausta ja errorien läpikäyntiä"
$ws.Range("D51").Value = "api"

$ws.Range("B52").Value = 1
$ws.Range("C52").Value = "dev ympäristön alustaminen, docker => konttiin backend + mongodb"
$ws.Range("D52").Value = "api"

$ws.Range("B53").Value = 2
$ws.Range("C53").Value = "dev ympäristön luotu loppuun, muutokset päivittyy suoraan konttiin, mongodb pyörii kontissa ongelmitta"
$ws.Range("D53").Value = "api"

# --- 3. Re-create the "tunnit yht." total row further down, at row 60 -----
$ws.Rows("60:60").RowHeight = 14.25
$ws.Range("A60").Value = "tunnit yht."
$ws.Range("B60").Formula = "=SUM(B2:B53)"

# --- 4. Selection / scroll position, matching where the user ended up -----
$excel.Goto($ws.Range("A43"), $true)
$ws.Range("B53").Select()

# --- 5. Page setup (print settings) ----------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
